$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTec logo, wp:docPr/name (and cNvPr/name) "image1.jpg" -> "image2.jpg" ---

# Default header (maps to header2.xml, docPr id="3")
$hDefault = $sec.Headers.Item(1).Range.InlineShapes.Item(1)
$hDefault.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

# First-page header (maps to header1.xml, docPr id="1")
$hFirst = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$hFirst.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

# --- Footers: Pearson logo, wp:docPr/name (and cNvPr/name) "image2.png" -> "image1.png" ---

# Default footer (maps to footer2.xml, docPr id="4")
$fDefault = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$fDefault.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# First-page footer (maps to footer1.xml, docPr id="2")
$fFirst = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$fFirst.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

Write-Host "Renamed header/footer logo images."
